$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Merge the two runs "TUE Feb 20" + " 13:11:14 PST 2018" into a
#    single run "TUE Feb 20 13:11:14 PST 2018".
# ---------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r") -eq "TUE Feb 20 13:11:14 PST 2018") {
        $s = $p.Range.Start
        $e = $p.Range.End
        $rr = $d.Range($s, $e)
        $rr.Text = "TUE Feb 20 13:11:14 PST 2018"
        break
    }
}

# ---------------------------------------------------------------
# 2) Insert the new "THU Feb 22" purchase record block right after
#    the "Amount balance ... - 167680.0" paragraph's first trailing
#    empty paragraph (i.e. between the existing empty PlainText
#    paragraph and the following empty bold PlainText paragraph).
# ---------------------------------------------------------------
$anchorIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r") -eq "Amount balance`t`t`t- 167680.0") {
        $anchorIndex = $i + 1
        break
    }
}

if ($anchorIndex -eq -1) {
    Write-Host "ERROR: anchor paragraph not found"
} else {
    # (isBold, text) pairs for the new paragraphs, in document order.
    $newParas = @(
        @($true,  ""),
        @($false, "THU Feb 22 13:03:00 PST 2018"),
        @($false, "Person Name`t`t`t`t- SY ABS"),
        @($false, "Bill number`t`t`t`t- 1306"),
        @($false, "---------------------------------------------------------------"),
        @($false, "Item Name`t`t`t`t- CARROT"),
        @($false, "Number of Pockets`t`t`t- 10"),
        @($false, "Number of KGs`t`t`t- 986"),
        @($false, "Rate`t`t`t`t`t- 9"),
        @($false, "Transport & Miscellaneous`t- 150"),
        @($false, "Total Price`t`t`t`t- 9024.0"),
        @($true,  "Amount balance`t`t`t- 176704.0"),
        @($false, "")
    )

    $curIndex = $anchorIndex
    foreach ($item in $newParas) {
        $isBold = $item[0]
        $text = $item[1]

        $anchor = $d.Paragraphs.Item($curIndex - 1).Range
        $anchor.Collapse(0)
        $anchor.InsertParagraphAfter()

        $newp = $d.Paragraphs.Item($curIndex)
        if ($isBold) {
            $newp.Range.Font.Bold = 1
        } else {
            $newp.Range.Font.Bold = 0
        }

        # Always type a placeholder char then replace/remove it, so
        # the paragraph never keeps a stray empty run, regardless of
        # whether it ends up with real text or stays blank.
        $newp.Range.InsertBefore("X")
        $s = $newp.Range.Start
        if ($text.Length -gt 0) {
            $e = $s + 1
            $rng = $d.Range($s, $e)
            $rng.Text = $text
        } else {
            $e = $s + 1
            $rng = $d.Range($s, $e)
            $rng.Delete()
        }

        $curIndex = $curIndex + 1
    }
}

Write-Host "Paragraphs now: $($d.Paragraphs.Count)"
